$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style (format) from the last existing data row (A255) down to the new date cells
$ws.Range("A255").Copy()
$ws.Range("A256:A269").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A256").Value = 44330
$ws.Range("B256").Value = 1
$ws.Range("C256").Value = 1
$ws.Range("D256").Value = 28.87669650591972

$ws.Range("A257").Value = 44331
$ws.Range("B257").Value = 1
$ws.Range("C257").Value = 2
$ws.Range("D257").Value = 57.75339301183945

$ws.Range("A258").Value = 44332
$ws.Range("B258").Value = 1
$ws.Range("C258").Value = 3
$ws.Range("D258").Value = 86.63008951775916

$ws.Range("A259").Value = 44333
$ws.Range("B259").Value = 1
$ws.Range("C259").Value = 4
$ws.Range("D259").Value = 115.5067860236789

$ws.Range("A260").Value = 44334
$ws.Range("B260").Value = 0
$ws.Range("C260").Value = 4
$ws.Range("D260").Value = 115.5067860236789

$ws.Range("A261").Value = 44335
$ws.Range("B261").Value = 0
$ws.Range("C261").Value = 4
$ws.Range("D261").Value = 115.5067860236789

$ws.Range("A262").Value = 44336
$ws.Range("B262").Value = 1
$ws.Range("C262").Value = 5
$ws.Range("D262").Value = 144.3834825295986

$ws.Range("A263").Value = 44337
$ws.Range("B263").Value = 0
$ws.Range("C263").Value = 4
$ws.Range("D263").Value = 115.5067860236789

$ws.Range("A264").Value = 44338
$ws.Range("B264").Value = 0
$ws.Range("C264").Value = 3
$ws.Range("D264").Value = 86.63008951775916

$ws.Range("A265").Value = 44339
$ws.Range("B265").Value = 0
$ws.Range("C265").Value = 2
$ws.Range("D265").Value = 57.75339301183945

$ws.Range("A266").Value = 44340
$ws.Range("B266").Value = 0
$ws.Range("C266").Value = 1
$ws.Range("D266").Value = 28.87669650591972

$ws.Range("A267").Value = 44341
$ws.Range("B267").Value = 0
$ws.Range("C267").Value = 1
$ws.Range("D267").Value = 28.87669650591972

$ws.Range("A268").Value = 44342
$ws.Range("B268").Value = 0
$ws.Range("C268").Value = 1
$ws.Range("D268").Value = 28.87669650591972

$ws.Range("A269").Value = 44343
$ws.Range("B269").Value = 0
$ws.Range("C269").Value = 0
$ws.Range("D269").Value = 0

Write-Output "done"
